$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the same cell formatting used by the existing "datetime" column
# (column A): bold Calibri font, centered alignment, thin box border and
# the workbook's custom date/time number format.
function Set-DateCell($addr, $dateVal) {
    $r = $ws.Range($addr)
    $r.Value = $dateVal
    $r.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $r.Font.Bold = $true
    $r.HorizontalAlignment = -4108
    $r.VerticalAlignment = -4160
    $r.Borders.LineStyle = 1
}

# Write one OHLCV-style data row: date in column A, the same numeric
# value repeated across open/high/low/close (B:E), and 0 volume (F).
function Set-DataRow($rowNum, $dateVal, $v) {
    Set-DateCell "A$rowNum" $dateVal
    $ws.Range("B$rowNum").Value = $v
    $ws.Range("C$rowNum").Value = $v
    $ws.Range("D$rowNum").Value = $v
    $ws.Range("E$rowNum").Value = $v
    $ws.Range("F$rowNum").Value = 0
}

# --- Step 1: append 3 new rows right after the current last data row (75) ---
# Doing this *before* the top insert means these rows inherit formatting
# from the row above them (still within the original data block), and
# they will be pushed down to rows 79-81 once the rows are inserted above.
$ws.Rows.Item(76).Insert()
$ws.Rows.Item(77).Insert()
$ws.Rows.Item(78).Insert()

Set-DataRow 76 45412 42842522000000
Set-DataRow 77 45443 42953758000000
Set-DataRow 78 45473 43654679000000

# --- Step 2: insert 3 new rows at the top of the second data block (row 39) ---
# This shifts the old rows 39-75 down to 42-78, and the rows appended in
# step 1 down to 79-81.
$ws.Rows.Item(39).Insert()
$ws.Rows.Item(40).Insert()
$ws.Rows.Item(41).Insert()

Set-DataRow 39 45412 42842522000000
Set-DataRow 40 45443 42953758000000
Set-DataRow 41 45473 43654679000000
